$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.514.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.923.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.61%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4843"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4100"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08190"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.025"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.59"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.902.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.058"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.259"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.17%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.008"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06774"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.53%  "
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.006"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.536.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.640"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.192"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.129.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.775"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.28%  "
$ws.Range("E28").Value = "  +2.64%  "
$ws.Range("E29").Value = "  +2.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.92%  "
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09591"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.539"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.570"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.87%  "
$ws.Range("E35").Value = "  +1.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02287"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06151"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.48%  "
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5996"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.064"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1866"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.42%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.284"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.21%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.407"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07615"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5603"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.963"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "117.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.441"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.98%  "
